# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet "2022-Q1" right before the "总计" (Total) sheet,
#    populated with the Q1-2022 fund holding snapshot (same layout as the
#    other quarterly sheets: 2020-Q4 / 2021-Q3 / 2021-Q4).
# 2. Update the "总计" (Total) summary sheet: add a new top row for
#    "2022-Q1" and push the previously existing rows down by one.

$wb = $excel.ActiveWorkbook

$templateSheet = $wb.Worksheets.Item("2021-Q4")

# --- 1. Create the new "2022-Q1" worksheet, positioned just before "总计" ---
$q1 = $wb.Worksheets.Add($wb.Worksheets.Item("总计"))
$q1.Name = "2022-Q1"

# NOTE: inserting a sheet shifts the position-based handle that
# "总计" used to resolve to, so re-fetch it by name *after* the insert
# to make sure later writes land on the right sheet.
$totalSheet = $wb.Worksheets.Item("总计")

# Copy the header formatting (bold + border + centered, style used by the
# other quarterly sheets) from the "2021-Q4" sheet onto the new sheet.
$templateSheet.Range("B1:H1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)

$templateSheet.Range("A2").Copy()
$q1.Range("A2").PasteSpecial(-4122)

# Header row
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# Data row (force the numeric-looking fund figures to stay text, matching
# how the other quarterly sheets store them)
$q1.Range("A2").Value = 0

$q1.Range("B2:G2").NumberFormat = "@"
$q1.Range("B2").Value = "513360"
$q1.Range("C2").Value = "博时中证全球中国教育主题交易型开放式指数证券投资基金(QDII)"
$q1.Range("D2").Value = "6.05"
$q1.Range("E2").Value = "99.49"
$q1.Range("F2").Value = "10.16"
$q1.Range("G2").Value = "0.6147"

$q1.Range("H2").Value = 1

# --- 2. Update the "总计" sheet with the new quarter on top ---

# Give the new A5 the same style as the existing index column cells before
# writing into it.
$totalSheet.Range("A4").Copy()
$totalSheet.Range("A5").PasteSpecial(-4122)

$totalSheet.Range("A5").Value = 3
$totalSheet.Range("B5").Value = "2020-Q4"
$totalSheet.Range("C5").Value = 7
$totalSheet.Range("D5").Value = 2.96

$totalSheet.Range("A4").Value = 2
$totalSheet.Range("B4").Value = "2021-Q3"
$totalSheet.Range("C4").Value = 1
$totalSheet.Range("D4").Value = 0.23

$totalSheet.Range("A3").Value = 1
$totalSheet.Range("B3").Value = "2021-Q4"
$totalSheet.Range("C3").Value = 1
$totalSheet.Range("D3").Value = 0.45

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 1
$totalSheet.Range("D2").Value = 0.61
